$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D5: append "(positive definite matrix)" to the title
$ws.Range("D5").Value = "양의 정부호 행렬 (positive definite matrix)"

# D9/E9: update title and link
$ws.Range("D9").Value = "해외 기업들 Data Scientist 공고 샘플 정리"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/data-scientist-job-postings/#utm_source=rss&utm_medium=rss&utm_campaign=data-scientist-job-postings"

# D26: update title
$ws.Range("D26").Value = "ai plus(est soft)"
